$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3457.439
$ws.Range("J17").Value = 3457.439
$ws.Range("L17").Value = 10372.317
$ws.Range("N17").Value = -10708.317
$ws.Range("H96").Value = 250028080
$ws.Range("I96").Value = 17266.5
$ws.Range("J96").Value = 500038880
$ws.Range("K96").Value = 51799.5
$ws.Range("L96").Value = 1500116640
$ws.Range("M96").Value = -50426.5
$ws.Range("N96").Value = -1500119386
$ws.Range("H138").Value = 1396.02
$ws.Range("I138").Value = 638.619
$ws.Range("J138").Value = 1944.4828
$ws.Range("K138").Value = 1915.857
$ws.Range("L138").Value = 5833.4484
$ws.Range("M138").Value = 3224.143
$ws.Range("N138").Value = -16113.4484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10646.351
$ws.Range("I32").Value = 9392.981
$ws.Range("K32").Value = 9392.981
$ws.Range("M32").Value = -9105.981
$ws.Range("H110").Value = 1771.7715
$ws.Range("I110").Value = 1746.2
$ws.Range("K110").Value = 1746.2
$ws.Range("M110").Value = 298.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3104.8096
$ws.Range("I105").Value = 2294.5454
$ws.Range("J105").Value = 3996.1
$ws.Range("K105").Value = 2294.5454
$ws.Range("L105").Value = 3996.1
$ws.Range("M105").Value = -547.5454
$ws.Range("N105").Value = -7490.1
$ws.Range("H107").Value = 1811.5278
$ws.Range("I107").Value = 1696.0358
$ws.Range("J107").Value = 2215.75
$ws.Range("K107").Value = 1696.0358
$ws.Range("L107").Value = 2215.75
$ws.Range("M107").Value = 223.9641999999999
$ws.Range("N107").Value = -6055.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3109.2
$ws.Range("I31").Value = 854.6842
$ws.Range("J31").Value = 3638.037
$ws.Range("K31").Value = 854.6842
$ws.Range("L31").Value = 3638.037
$ws.Range("M31").Value = -559.6842
$ws.Range("N31").Value = -4228.037
$ws.Range("H34").Value = 3109.2
$ws.Range("I34").Value = 854.6842
$ws.Range("J34").Value = 3638.037
$ws.Range("K34").Value = 854.6842
$ws.Range("L34").Value = 3638.037
$ws.Range("M34").Value = -652.6842
$ws.Range("N34").Value = -4042.037
$ws.Range("H58").Value = 1471.3903
$ws.Range("I58").Value = 1144.6207
$ws.Range("J58").Value = 2261.0833
$ws.Range("K58").Value = 1144.6207
$ws.Range("L58").Value = 2261.0833
$ws.Range("M58").Value = -941.6206999999999
$ws.Range("N58").Value = -2667.0833
$ws.Range("H132").Value = 40828.168
$ws.Range("I132").Value = 1580.4231
$ws.Range("J132").Value = 142872.3
$ws.Range("K132").Value = 4741.2693
$ws.Range("L132").Value = 428616.9
$ws.Range("M132").Value = -2211.2693
$ws.Range("N132").Value = -433676.9
$ws.Range("H134").Value = 484246.4
$ws.Range("I134").Value = 1360.4286
$ws.Range("J134").Value = 1751822.1
$ws.Range("K134").Value = 4081.2858
$ws.Range("L134").Value = 5255466.300000001
$ws.Range("M134").Value = -1546.2858
$ws.Range("N134").Value = -5260536.300000001
$ws.Range("H136").Value = 1471.3903
$ws.Range("I136").Value = 1144.6207
$ws.Range("J136").Value = 2261.0833
$ws.Range("K136").Value = 3433.8621
$ws.Range("L136").Value = 6783.249899999999
$ws.Range("M136").Value = -883.8620999999998
$ws.Range("N136").Value = -11883.2499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5675.3335
$ws.Range("I5").Value = 7233.1333
$ws.Range("K5").Value = 21699.3999
$ws.Range("M5").Value = -21587.3999
$ws.Range("H122").Value = 8010.6665
$ws.Range("J122").Value = 16509.285
$ws.Range("L122").Value = 148583.565
$ws.Range("N122").Value = -153483.565
$ws.Range("H131").Value = 3768.3076
$ws.Range("I131").Value = 9475.637000000001
$ws.Range("J131").Value = 1526.1428
$ws.Range("K131").Value = 28426.911
$ws.Range("L131").Value = 4578.428400000001
$ws.Range("M131").Value = -23386.911
$ws.Range("N131").Value = -14658.4284
$ws.Range("H135").Value = 5675.3335
$ws.Range("I135").Value = 7233.1333
$ws.Range("K135").Value = 65098.1997
$ws.Range("M135").Value = -62563.1997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4813.793
$ws.Range("I70").Value = 4785.185
$ws.Range("J70").Value = 5200
$ws.Range("K70").Value = 4785.185
$ws.Range("L70").Value = 5200
$ws.Range("M70").Value = -4515.185
$ws.Range("N70").Value = -5740
$ws.Range("H73").Value = 4813.793
$ws.Range("I73").Value = 4785.185
$ws.Range("J73").Value = 5200
$ws.Range("K73").Value = 4785.185
$ws.Range("L73").Value = 5200
$ws.Range("M73").Value = -3849.185
$ws.Range("N73").Value = -7072
$ws.Range("H132").Value = 2611.4546
$ws.Range("I132").Value = 1574.8235
$ws.Range("K132").Value = 4724.470499999999
$ws.Range("M132").Value = -2194.470499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2676.8
$ws.Range("I61").Value = 2718.111
$ws.Range("J61").Value = 2305
$ws.Range("K61").Value = 2718.111
$ws.Range("L61").Value = 2305
$ws.Range("M61").Value = -2516.111
$ws.Range("N61").Value = -2709
$ws.Range("H113").Value = 2676.8
$ws.Range("I113").Value = 2718.111
$ws.Range("J113").Value = 2305
$ws.Range("K113").Value = 2718.111
$ws.Range("L113").Value = 2305
$ws.Range("M113").Value = -548.1109999999999
$ws.Range("N113").Value = -6645
$ws.Range("H122").Value = 64695.812
$ws.Range("I122").Value = 85352.75
$ws.Range("J122").Value = 2725
$ws.Range("K122").Value = 256058.25
$ws.Range("L122").Value = 8175
$ws.Range("M122").Value = -253608.25
$ws.Range("N122").Value = -13075
$ws.Range("H124").Value = 46993.668
$ws.Range("J124").Value = 46993.668
$ws.Range("L124").Value = 46993.668
$ws.Range("N124").Value = -56813.668
$ws.Range("H136").Value = 1652.909
$ws.Range("I136").Value = 1335.9565
$ws.Range("K136").Value = 4007.8695
$ws.Range("M136").Value = -1457.8695

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 595.53845
$ws.Range("I113").Value = 541.55554
$ws.Range("J113").Value = 717
$ws.Range("K113").Value = 1624.66662
$ws.Range("L113").Value = 2151
$ws.Range("M113").Value = 545.33338
$ws.Range("N113").Value = -6491
$ws.Range("H138").Value = 43385.715
$ws.Range("J138").Value = 43385.715
$ws.Range("L138").Value = 43385.715
$ws.Range("N138").Value = -53665.715
